$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Clear the previously-computed address-parsing results (Building Name, Unit,
# Street No, Street Name, Street Type) for the existing rows 2-16.
# Street No (column H) is fully removed, while Building Name (F), Unit (G),
# Street Name (I) and Street Type (J) remain present but empty.
$ws.Range("F2:J16").ClearContents()
$ws.Range("F2:F16").Font.Bold = $false
$ws.Range("G2:G16").Font.Bold = $false
$ws.Range("I2:I16").Font.Bold = $false
$ws.Range("J2:J16").Font.Bold = $false

# Add two new customer rows.
$ws.Range("A17").Value = 126
$ws.Range("A17").Style = "Normal"
$ws.Range("B17").Value = "Customer 16"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "Unit 31 A 53 BCH Street"
$ws.Range("C17").Style = "Normal"
$ws.Range("E17").Value = "Wollongong"
$ws.Range("E17").Style = "Normal"
$ws.Range("F17").Font.Bold = $false
$ws.Range("G17").Font.Bold = $false
$ws.Range("I17").Font.Bold = $false
$ws.Range("J17").Font.Bold = $false

$ws.Range("A18").Value = 127
$ws.Range("A18").Style = "Normal"
$ws.Range("B18").Value = "Customer 17"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "Unit 11 Abert Building 52 Abert Street"
$ws.Range("C18").Style = "Normal"
$ws.Range("F18").Font.Bold = $false
$ws.Range("G18").Font.Bold = $false
$ws.Range("I18").Font.Bold = $false
$ws.Range("J18").Font.Bold = $false

# Update the selection to match the author's saved state.
$ws.Range("F2:J18").Select()
